$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for data rows 2-61 changes from serial date 45179 to 45180
for ($r = 2; $r -le 61; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
